$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C text updates (coin re-ranking) ---
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

# --- Column D price updates ---
$ws.Range("D2").Value = "26.641.69"
$ws.Range("D3").Value = "1.594.76"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.22"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0618"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.68"
$ws.Range("D12").Value = "1.814.53"
$ws.Range("D13").Value = "1.593.16"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.525"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.13"
$ws.Range("D17").Value = "26.601.06"
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.05"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.00"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.30"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.90"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.21"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.15"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.36"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0506"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.666"
$ws.Range("D35").Value = "1.302.81"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.833"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.792"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.36"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.51"
$ws.Range("D45").Value = "1.727.58"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.32"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.821"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0988"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0503"

# --- Column E volume(1h) % updates ---
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("E23").Value = "  -3.82%  "
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("E33").Value = "  -10.34%  "
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  -4.64%  "
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("E51").Value = "  -2.61%  "
